$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HR")
[void]$ws.Select()

# New "Movement Registry" sub-menu row under HR > Management.
# Write the URL first, then the label, so the new shared-string entries
# land in the same order the workbook author produced them in.
$ws.Range("D66").Value = "/hr/movement-registry/"
$ws.Range("A66").Value = "Movement Registry"
$ws.Range("B66").Value = "Yes"
$ws.Range("C66").Value = "Management"
$ws.Range("E66").Value = "fas fa-users"
$ws.Range("F66").Value = 27

# Scroll the HR sheet so row 52 is at the top and select E66 (the new row),
# matching where the author ended up after adding the entry.
$w = $excel.ActiveWindow
$w.ScrollRow = 52
$w.ScrollColumn = 1
[void]$ws.Range("E66").Select()
